# Daily attendance processing - swap "Recorded By" order for rows that list
# both the staff email and "System", turning "email, System" into
# "System, email" (column G - "Recorded By").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$colG = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
